# Legislator property workbook - #5: property boat&car done
# Rebuild the "汽車" (car) sheet (sheet3): turn row 1 from a (wrong) data
# duplicate into the standard header row, and append the standard
# metadata columns (H:N) to the two data rows, matching the other
# property sheets (land/building/...) in this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# ---- Row 1: standard header labels (B1:N1) -------------------------------
$ws.Cells.Item(1, 2).Value  = "name"
$ws.Cells.Item(1, 3).Value  = "capacity"
$ws.Cells.Item(1, 4).Value  = "owner"
$ws.Cells.Item(1, 5).Value  = "register_date"
$ws.Cells.Item(1, 6).Value  = "register_reason"
$ws.Cells.Item(1, 7).Value  = "acquire_value"
$ws.Cells.Item(1, 8).Value  = "property_category"
$ws.Cells.Item(1, 9).Value  = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Give the new header cells (H1:N1) the same look as the existing header
# cells (bold font + border) by copying the format from an existing
# header cell instead of guessing at style indices.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("H1:N1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---- Row 2: car #48 (Volvo) - fill in the metadata columns ---------------
$ws.Cells.Item(2, 8).Value  = "land"
$ws.Cells.Item(2, 9).Value  = "normal"
$ws.Cells.Item(2, 10).Value = "2012-04-18"
$ws.Cells.Item(2, 11).Value = "李應元"
$ws.Cells.Item(2, 12).Value = 708
$ws.Cells.Item(2, 13).Value = "tmp3fed1"
$ws.Cells.Item(2, 14).Value = 48

# ---- Row 3: car #49 (Toyota Camary) - fill in the metadata columns -------
$ws.Cells.Item(3, 8).Value  = "land"
$ws.Cells.Item(3, 9).Value  = "normal"
$ws.Cells.Item(3, 10).Value = "2012-04-18"
$ws.Cells.Item(3, 11).Value = "李應元"
$ws.Cells.Item(3, 12).Value = 708
$ws.Cells.Item(3, 13).Value = "tmp3fed1"
$ws.Cells.Item(3, 14).Value = 49

# Give the new data cells (H2:N3) the same look as the existing data
# cells in that row.
$ws.Range("G2").Copy() | Out-Null
$ws.Range("H2:N2").PasteSpecial(-4122) | Out-Null
$ws.Range("G3").Copy() | Out-Null
$ws.Range("H3:N3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
